$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 12:04"

# Update the province data table values
$ws.Range("B4").Value = 63870
$ws.Range("C4").Value = 38670
$ws.Range("D4").Value = 16696
$ws.Range("E4").Value = 8504

$ws.Range("B5").Value = 51190
$ws.Range("C5").Value = 23231
$ws.Range("D5").Value = 22565
$ws.Range("E5").Value = 5394

$ws.Range("B6").Value = 17625
$ws.Range("C6").Value = 7149
$ws.Range("D6").Value = 8612
$ws.Range("E6").Value = 1864

$ws.Range("B7").Value = 16184
$ws.Range("C7").Value = 5895
$ws.Range("D7").Value = 7612
$ws.Range("E7").Value = 2677

$ws.Range("B9").Value = 12268
$ws.Range("C9").Value = 8076
$ws.Range("D9").Value = 2898
$ws.Range("E9").Value = 1294

$ws.Range("B10").Value = 9134
$ws.Range("C10").Value = 6959
$ws.Range("D10").Value = 1589
$ws.Range("E10").Value = 586

$ws.Range("B13").Value = 5258
$ws.Range("C13").Value = 2880
$ws.Range("D13").Value = 1578
$ws.Range("E13").Value = 800

$ws.Range("B15").Value = 4983
$ws.Range("C15").Value = 2732
$ws.Range("D15").Value = 1771
$ws.Range("E15").Value = 480

$ws.Range("B16").Value = 3986
$ws.Range("C16").Value = 2457
$ws.Range("D16").Value = 1191
$ws.Range("E16").Value = 338

$ws.Range("B23").Value = 2877
$ws.Range("C23").Value = 2265
$ws.Range("D23").Value = 145
$ws.Range("E23").Value = 467

$ws.Range("B30").Value = 2326
$ws.Range("C30").Value = 971
$ws.Range("D30").Value = 1063
$ws.Range("E30").Value = 292

$ws.Range("C31").Value = 1258
$ws.Range("D31").Value = 833

$ws.Range("B33").Value = 2220
$ws.Range("C33").Value = 1758
$ws.Range("D33").Value = 262
$ws.Range("E33").Value = 200

$ws.Range("C59").Value = 108
$ws.Range("D59").Value = 9
